$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4910.6523
$ws.Range("I62").Value = 4797.4
$ws.Range("J62").Value = 5123
$ws.Range("K62").Value = 4797.4
$ws.Range("L62").Value = 5123
$ws.Range("M62").Value = -4173.4
$ws.Range("N62").Value = -6371

$ws.Range("H65").Value = 4910.6523
$ws.Range("I65").Value = 4797.4
$ws.Range("J65").Value = 5123
$ws.Range("K65").Value = 23987
$ws.Range("L65").Value = 25615
$ws.Range("M65").Value = -20867
$ws.Range("N65").Value = -31855

$ws.Range("H74").Value = 5173.5
$ws.Range("I74").Value = 6678.625
$ws.Range("J74").Value = 3166.6667
$ws.Range("K74").Value = 6678.625
$ws.Range("L74").Value = 3166.6667
$ws.Range("M74").Value = -5742.625
$ws.Range("N74").Value = -5038.6667

$ws.Range("H77").Value = 5173.5
$ws.Range("I77").Value = 6678.625
$ws.Range("J77").Value = 3166.6667
$ws.Range("K77").Value = 33393.125
$ws.Range("L77").Value = 15833.3335
$ws.Range("M77").Value = -28713.125
$ws.Range("N77").Value = -25193.3335

$ws.Range("H140").Value = 99984
$ws.Range("J140").Value = 99984
$ws.Range("L140").Value = 99984
$ws.Range("N140").Value = -110344

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4933.7188
$ws.Range("I61").Value = 1286.08
$ws.Range("J61").Value = 17961
$ws.Range("K61").Value = 1286.08
$ws.Range("L61").Value = 17961
$ws.Range("M61").Value = -1074.08
$ws.Range("N61").Value = -18385

$ws.Range("H97").Value = 750.8378
$ws.Range("I97").Value = 700.34375
$ws.Range("J97").Value = 1074
$ws.Range("K97").Value = 700.34375
$ws.Range("L97").Value = 1074
$ws.Range("M97").Value = -204.34375
$ws.Range("N97").Value = -2066

$ws.Range("H122").Value = 2996.2856
$ws.Range("I122").Value = 2829
$ws.Range("K122").Value = 8487
$ws.Range("M122").Value = -6037

$ws.Range("H132").Value = 2453.077
$ws.Range("I132").Value = 1432.3334
$ws.Range("K132").Value = 4297.0002
$ws.Range("M132").Value = -1767.0002

$ws.Range("H136").Value = 4933.7188
$ws.Range("I136").Value = 1286.08
$ws.Range("J136").Value = 17961
$ws.Range("K136").Value = 3858.24
$ws.Range("L136").Value = 53883
$ws.Range("M136").Value = -1308.24
$ws.Range("N136").Value = -58983

$ws.Range("H138").Value = 94194.25
$ws.Range("J138").Value = 94194.25
$ws.Range("L138").Value = 94194.25
$ws.Range("N138").Value = -104474.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6732.3716
$ws.Range("I20").Value = 7806.9653
$ws.Range("K20").Value = 7806.9653
$ws.Range("M20").Value = -7559.9653

$ws.Range("H134").Value = 3274.5
$ws.Range("I134").Value = 3035.75
$ws.Range("J134").Value = 4707
$ws.Range("K134").Value = 9107.25
$ws.Range("L134").Value = 14121
$ws.Range("M134").Value = -6572.25
$ws.Range("N134").Value = -19191

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 885.4
$ws.Range("I5").Value = 425
$ws.Range("K5").Value = 425
$ws.Range("M5").Value = -313

$ws.Range("H8").Value = 6999.5
$ws.Range("I8").Value = 6999.5
$ws.Range("K8").Value = 6999.5
$ws.Range("M8").Value = -6859.5

$ws.Range("H31").Value = 14292628
$ws.Range("I31").Value = 20006980
$ws.Range("K31").Value = 20006980
$ws.Range("M31").Value = -20006685

$ws.Range("H34").Value = 14292628
$ws.Range("I34").Value = 20006980
$ws.Range("K34").Value = 20006980
$ws.Range("M34").Value = -20006778

$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 5000
$ws.Range("K39").Value = 5000
$ws.Range("M39").Value = -4609

$ws.Range("H41").Value = 40000
$ws.Range("J41").Value = 40000
$ws.Range("L41").Value = 40000
$ws.Range("N41").Value = -40856

$ws.Range("H49").Value = 5000
$ws.Range("I49").Value = 5000
$ws.Range("K49").Value = 5000
$ws.Range("M49").Value = -4818

$ws.Range("H86").Value = 49561.812
$ws.Range("I86").Value = 64090.363
$ws.Range("K86").Value = 64090.363
$ws.Range("M86").Value = -62967.363

$ws.Range("H89").Value = 49561.812
$ws.Range("I89").Value = 64090.363
$ws.Range("K89").Value = 320451.815
$ws.Range("M89").Value = -314835.815

$ws.Range("H99").Value = 2928.9167
$ws.Range("I99").Value = 2171.1428
$ws.Range("K99").Value = 2171.1428
$ws.Range("M99").Value = -673.1428000000001

$ws.Range("H126").Value = 2928.9167
$ws.Range("I126").Value = 2171.1428
$ws.Range("K126").Value = 6513.428400000001
$ws.Range("M126").Value = -4043.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 4999.75
$ws.Range("J88").Value = 4999.75
$ws.Range("L88").Value = 14999.25
$ws.Range("N88").Value = -15855.25

$ws.Range("H91").Value = 4999.75
$ws.Range("J91").Value = 4999.75
$ws.Range("L91").Value = 14999.25
$ws.Range("N91").Value = -17963.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 17818.35
$ws.Range("J24").Value = 18177.264
$ws.Range("L24").Value = 18177.264
$ws.Range("N24").Value = -18523.264

$ws.Range("H29").Value = 20666
$ws.Range("J29").Value = 20666
$ws.Range("L29").Value = 20666
$ws.Range("N29").Value = -21246

$ws.Range("H126").Value = 2636.4285
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 3966.5557
$ws.Range("I132").Value = 3339
$ws.Range("K132").Value = 10017
$ws.Range("M132").Value = -7487

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 43333.332
$ws.Range("I4").Value = 25000
$ws.Range("K4").Value = 25000
$ws.Range("M4").Value = -24887

$ws.Range("H14").Value = 10960.125
$ws.Range("I14").Value = 8416.75
$ws.Range("J14").Value = 13503.5
$ws.Range("K14").Value = 8416.75
$ws.Range("L14").Value = 13503.5
$ws.Range("M14").Value = -8244.75
$ws.Range("N14").Value = -13847.5

$ws.Range("H24").Value = 13537.4
$ws.Range("I24").Value = 11996.75
$ws.Range("K24").Value = 11996.75
$ws.Range("M24").Value = -11653.75

$ws.Range("H26").Value = 100000
$ws.Range("J26").Value = 100000
$ws.Range("L26").Value = 100000
$ws.Range("N26").Value = -100590

$ws.Range("H28").Value = 43333.332
$ws.Range("I28").Value = 25000
$ws.Range("K28").Value = 25000
$ws.Range("M28").Value = -24768

$ws.Range("H31").Value = 3496.5
$ws.Range("I31").Value = 5000
$ws.Range("K31").Value = 5000
$ws.Range("M31").Value = -4752

$ws.Range("H37").Value = 43333.332
$ws.Range("I37").Value = 25000
$ws.Range("K37").Value = 25000
$ws.Range("M37").Value = -24893

$ws.Range("H46").Value = 3747.2856
$ws.Range("I46").Value = 2167
$ws.Range("J46").Value = 4178.273
$ws.Range("K46").Value = 2167
$ws.Range("L46").Value = 4178.273
$ws.Range("M46").Value = -1979
$ws.Range("N46").Value = -4554.273

$ws.Range("H100").Value = 2343.5715
$ws.Range("I100").Value = 2182.7273
$ws.Range("K100").Value = 2182.7273
$ws.Range("M100").Value = -1641.7273

$ws.Range("H122").Value = 2849.9583
$ws.Range("I122").Value = 2895.2856
$ws.Range("J122").Value = 2532.6667
$ws.Range("K122").Value = 8685.856800000001
$ws.Range("L122").Value = 7598.000100000001
$ws.Range("M122").Value = -6235.856800000001
$ws.Range("N122").Value = -12498.0001

$ws.Range("H132").Value = 5014.8945
$ws.Range("I132").Value = 4791.077
$ws.Range("K132").Value = 14373.231
$ws.Range("M132").Value = -11843.231

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H136").Value = 5312.722
$ws.Range("J136").Value = 5593.625
$ws.Range("L136").Value = 16780.875
$ws.Range("N136").Value = -21880.875

$ws.Range("H137").Value = 120000
$ws.Range("J137").Value = 120000
$ws.Range("L137").Value = 120000
$ws.Range("N137").Value = -130200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 21534.5
$ws.Range("J18").Value = 23069
$ws.Range("L18").Value = 23069
$ws.Range("N18").Value = -23415

$ws.Range("H122").Value = 85441.42999999999
$ws.Range("I122").Value = 91895.58
$ws.Range("J122").Value = 1537.5
$ws.Range("K122").Value = 275686.74
$ws.Range("L122").Value = 4612.5
$ws.Range("M122").Value = -273236.74
$ws.Range("N122").Value = -9512.5

$ws.Range("H132").Value = 25755.139
$ws.Range("I132").Value = 27329.592
$ws.Range("K132").Value = 81988.776
$ws.Range("M132").Value = -79458.776

$ws.Range("H136").Value = 34502.094
$ws.Range("I136").Value = 39674.668
$ws.Range("J136").Value = 3466.6667
$ws.Range("K136").Value = 119024.004
$ws.Range("L136").Value = 10400.0001
$ws.Range("M136").Value = -116474.004
$ws.Range("N136").Value = -15500.0001
